$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell "D2" "28.309.99"
Set-TextCell "E2" "  -0.47%  "
Set-TextCell "D3" "1.577.25"
Set-TextCell "E3" "  +0.18%  "
Set-TextCell "E4" "  +0.30%  "
Set-TextCell "D5" "212.32"
Set-TextCell "E5" "  +0.79%  "
Set-TextCell "D6" "0.489"
Set-TextCell "E6" "  -0.20%  "
Set-TextCell "E7" "  +0.36%  "
Set-TextCell "D8" "44.49"
Set-TextCell "E8" "  -3.76%  "
Set-TextCell "D9" "23.81"
Set-TextCell "E9" "  +0.26%  "
Set-TextCell "E10" "  -0.55%  "
Set-TextCell "E11" "  -0.38%  "
Set-TextCell "D12" "0.0895"
Set-TextCell "E12" "  +1.81%  "
Set-TextCell "D13" "1.803.64"
Set-TextCell "E13" "  +0.20%  "
Set-TextCell "D14" "1.575.43"
Set-TextCell "E14" "  +0.01%  "
Set-TextCell "D15" "3.68"
Set-TextCell "E15" "  -0.53%  "
Set-TextCell "E16" "  -0.76%  "
Set-TextCell "D17" "28.350.94"
Set-TextCell "E17" "  -0.36%  "
Set-TextCell "D18" "61.71"
Set-TextCell "E18" "  -0.95%  "
Set-TextCell "D19" "230.89"
Set-TextCell "E19" "  +0.94%  "
Set-TextCell "E20" "  +1.17%  "
Set-TextCell "D21" "0.0₃0686"
Set-TextCell "E22" "  +0.38%  "
Set-TextCell "E23" "  +0.36%  "
Set-TextCell "D24" "9.05"
Set-TextCell "E24" "  -1.07%  "
Set-TextCell "E25" "  +3.19%  "
Set-TextCell "D26" "151.92"
Set-TextCell "E26" "  +0.63%  "
Set-TextCell "D27" "15.00"
Set-TextCell "E27" "  -0.12%  "
Set-TextCell "E28" "  -1.10%  "
Set-TextCell "E29" "  -0.66%  "
Set-TextCell "E30" "  +0.30%  "
Set-TextCell "D31" "0.0480"
Set-TextCell "E31" "  +3.54%  "
Set-TextCell "E32" "  -2.87%  "
Set-TextCell "E33" "  -0.03%  "
Set-TextCell "E34" "  -0.99%  "
Set-TextCell "D35" "1.389.93"
Set-TextCell "E35" "  -0.09%  "
Set-TextCell "E36" "  +6.90%  "
Set-TextCell "E37" "  -2.46%  "
Set-TextCell "E38" "  +0.21%  "
Set-TextCell "D39" "2.64"
Set-TextCell "E39" "  +3.50%  "
Set-TextCell "E40" "  -1.14%  "
Set-TextCell "D41" "0.520"
Set-TextCell "E41" "  -2.14%  "
Set-TextCell "E42" "  +0.35%  "
Set-TextCell "D43" "1.89"
Set-TextCell "E43" "  +2.95%  "
Set-TextCell "E44" "  -0.71%  "
Set-TextCell "D45" "0.0457"
Set-TextCell "E45" "  -2.33%  "
Set-TextCell "D46" "5.40"
Set-TextCell "E46" "  -3.83%  "
Set-TextCell "D47" "0.926"
Set-TextCell "E47" "  -5.33%  "
Set-TextCell "D48" "62.45"
Set-TextCell "E48" "  +0.17%  "
Set-TextCell "D49" "1.715.19"
Set-TextCell "E49" "  +0.23%  "
Set-TextCell "D50" "85.44"
Set-TextCell "E50" "  -0.25%  "
Set-TextCell "B51" "BabyDogeCoin"
Set-TextCell "C51" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextCell "D51" "0.0₆0102"
Set-TextCell "E51" "  -0.35%  "
